$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete one-year-ahead forecast values for the very first
# two rows (the naive component forecaster previously emitted bogus values
# here before the bug fix).
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Corrected forecast values produced after fixing the naive component
# forecaster bug (tiny floating point corrections throughout the series).
$ws.Range("E3").Value = 2.64925711235009
$ws.Range("C4").Value = 2.533533936850585
$ws.Range("E5").Value = 2.332261646026246
$ws.Range("C6").Value = 1.21254482274098
$ws.Range("E6").Value = 1.839804681163337
$ws.Range("E7").Value = 0.6705904529405782
$ws.Range("C8").Value = 0.4712609263772816
$ws.Range("E8").Value = 0.8520644823059031
$ws.Range("C11").Value = 4.109890522944326
$ws.Range("E11").Value = 3.628019428949014
$ws.Range("E15").Value = 3.933586883651397
$ws.Range("C16").Value = 2.777797690741446
$ws.Range("E16").Value = 2.073300717643911
$ws.Range("E17").Value = 1.589741018019186
$ws.Range("C18").Value = -1.432689847121826
$ws.Range("C19").Value = 2.033479419175155
$ws.Range("E19").Value = 1.562315774899048
